$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: some Price values look numeric (e.g. "19.41") and would be auto-
# converted to a number by Excel unless the cell is explicitly formatted as
# Text first - matches the source data which stores them as plain strings.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.298.52"
$ws.Range("E2").Value = "  -0.60%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.588.34"
$ws.Range("E3").Value = "  -0.26%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.64%  "

# Row 5 - BNB
$ws.Range("E5").Value = "  +0.10%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.05%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.58%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  +0.60%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -0.29%  "

# Row 10 - Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.41"
$ws.Range("E10").Value = "  -0.41%  "

# Row 11 - TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0847"
$ws.Range("E11").Value = "  +0.63%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("E12").Value = "  -0.18%  "

# Row 13 / 14 - Polkadot and WrappedEther swap places (with new price/volume data)
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.601.12"
$ws.Range("E13").Value = "  +0.49%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.06"
$ws.Range("E14").Value = "  +1.14%  "

# Row 15 - Polygon
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.520"
$ws.Range("E15").Value = "  +0.65%  "

# Row 16 - Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.33"
$ws.Range("E16").Value = "  +0.38%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "26.309.40"

# Row 18 - ShibaInu
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("E18").Value = "  -0.76%  "

# Row 19 - Chainlink
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.47"
$ws.Range("E19").Value = "  +6.23%  "

# Row 20 - BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "210.79"
$ws.Range("E20").Value = "  +1.99%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.66%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +0.31%  "

# Row 23 - Avalanche
$ws.Range("E23").Value = "  +1.15%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -2.55%  "

# Row 25 - Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.47"
$ws.Range("E25").Value = "  -0.07%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  -0.64%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  +0.23%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  +0.13%  "

# Row 29 - EthereumClassic
$ws.Range("E29").Value = "  +0.25%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +0.58%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +0.16%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -0.63%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +2.07%  "

# Row 34 - Maker
$ws.Range("D34").Value = "1.317.06"
$ws.Range("E34").Value = "  +2.88%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  -1.71%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  +2.06%  "

# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  -0.09%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  +0.85%  "

# Row 39 - WEMIXToken
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.09"
$ws.Range("E39").Value = "  -12.16%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  -1.01%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  -0.64%  "

# Row 42 - FraxShare
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.64"
$ws.Range("E42").Value = "  +4.52%  "

# Row 43 - TrustWalletToken
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.767"
$ws.Range("E43").Value = "  -0.13%  "

# Row 44 - MXToken
$ws.Range("E44").Value = "  -0.51%  "

# Row 45 - Aave
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.27"
$ws.Range("E45").Value = "  +0.09%  "

# Row 46 - RocketPoolETH
$ws.Range("D46").Value = "1.725.02"
$ws.Range("E46").Value = "  -0.07%  "

# Row 47 - Quant
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.36"
$ws.Range("E47").Value = "  -1.46%  "

# Row 48 - RenderToken
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.48"
$ws.Range("E48").Value = "  -4.64%  "

# Row 49 - Cronos
$ws.Range("E49").Value = "  -1.29%  "

# Row 50 - Algorand
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0976"
$ws.Range("E50").Value = "  -3.97%  "

# Row 51 - USDD
$ws.Range("E51").Value = "  -0.77%  "
